$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: create a brand-new cell with a specific look by cloning the number
# format / font / fill / border (but not contents) from an existing cell that
# already carries the desired style, then writing the new value into it.
# ---------------------------------------------------------------------------
function Set-StyledValue {
    param(
        [string]$SourceAddress,
        [string]$TargetAddress,
        $Value
    )
    $ws.Range($SourceAddress).Copy() | Out-Null
    $ws.Range($TargetAddress).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    $ws.Range($TargetAddress).Value = $Value
}
# NOTE: this engine's PowerShell subset does not bind named (-Param value)
# arguments on user-defined functions, only positional ones -- so every call
# below passes SourceAddress/TargetAddress/Value positionally.

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Cells that already exist (just empty) -- only the value needs to be filled
# in, the cell keeps its current style.
# ---------------------------------------------------------------------------
$ws.Range("D5").Value  = 5
$ws.Range("K5").Value  = 3
$ws.Range("K6").Value  = 3
$ws.Range("K17").Value = 1
$ws.Range("K20").Value = 4
$ws.Range("I21").Value = 5
$ws.Range("K21").Value = 2
$ws.Range("K22").Value = 3
$ws.Range("E24").Value = 5
$ws.Range("K26").Value = 4
$ws.Range("K27").Value = 4
$ws.Range("K28").Value = 4
$ws.Range("K33").Value = 1

# ---------------------------------------------------------------------------
# Brand-new cells -- style needs to be cloned from a neighbouring cell that
# already has the right look, then the value is written.
# ---------------------------------------------------------------------------
Set-StyledValue "I12" "I5"  5
Set-StyledValue "I12" "J6"  5
Set-StyledValue "I12" "J17" 5
Set-StyledValue "I12" "J22" 5
Set-StyledValue "I12" "H24" 5
Set-StyledValue "I29" "J24" 5
Set-StyledValue "I12" "J26" 5
Set-StyledValue "I12" "J27" 5
Set-StyledValue "I12" "J28" 5
Set-StyledValue "G29" "I33" 5

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# View state: move the active selection to J21 (the frozen-pane scroll
# position follows the engine's own bookkeeping).
# ---------------------------------------------------------------------------
$ws.Range("J21").Select()
